$wb = $excel.ActiveWorkbook

# Add the two new sheets at the end of the workbook, preserving order:
# CaseDetailStat, then CaseDetailStat_Message
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$caseDetailStat = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$caseDetailStat.Name = "CaseDetailStat"

$caseDetailStatMessage = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $caseDetailStat)
$caseDetailStatMessage.Name = "CaseDetailStat_Message"

# CaseDetailStat stays empty (just a blank row 1) - nothing to write.

# CaseDetailStat_Message: log of the (failed) CaseDetailStat cypher run.
# Block 1: Neo4j connection info + the CaseDetailStat query + output path.
$trialsCypher = "MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report) WITH DISTINCT c AS c, t ,a, s WHERE c.disease IN ['Adenocarcinoma of the cervix'] RETURN coalesce(c.case_id,'') AS ``Case ID`` , coalesce(t.clinical_trial_designation ,'')as ``Trial Code`` , coalesce(a.arm_id,'') As ``Arm`` , coalesce(a.arm_drug,'') As ``Arm Treatment`` , coalesce(c.disease,'') As Diagnosis , coalesce(c.gender,'') As Gender , coalesce(c.race,'') As Race , coalesce(c.ethnicity,'') As Ethnicity"
$statCypher = "MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report) WITH DISTINCT c AS c, t ,a, s WHERE c.disease IN ['Adenocarcinoma of the cervix'] OPTIONAL MATCH (s)<-[*]-(f:file) RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(t.clinical_trial_designation)) as number_of_trial"
$outputPath = "C:\Users\radhakrishnang2\Desktop\Commons_Automation\OutputFiles\TC01_Trials_Filter_Diagnosis-AdenoCervix_Neo4jData.xlsx"

$ws = $caseDetailStatMessage

$ws.Range("A1").Value = "Neo4j_URL:"
$ws.Range("A2").Value = "bolt://ncidb-q325-c.nci.nih.gov:7687"
$ws.Range("A3").Value = "User_name:"
$ws.Range("A4").Value = "neo4j"
$ws.Range("A5").Value = "PWD:"
$ws.Range("A6").Value = "icdcDBneo4j0"
$ws.Range("A7").Value = "Cypher:"
$ws.Range("A8").Value = $trialsCypher
$ws.Range("A9").Value = "Output:"
$ws.Range("A10").Value = $outputPath

# Block 2: same connection info + the stat cypher + output path.
$ws.Range("A11").Value = "Neo4j_URL:"
$ws.Range("A12").Value = "bolt://ncidb-q325-c.nci.nih.gov:7687"
$ws.Range("A13").Value = "User_name:"
$ws.Range("A14").Value = "neo4j"
$ws.Range("A15").Value = "PWD:"
$ws.Range("A16").Value = "icdcDBneo4j0"
$ws.Range("A17").Value = "Cypher:"
$ws.Range("A18").Value = $statCypher
$ws.Range("A19").Value = "Output:"
$ws.Range("A20").Value = $outputPath

# Block 3: connection info, error (empty cypher), output path.
$ws.Range("A21").Value = "Cypher query should not be an empty string"
$ws.Range("A22").Value = "Neo4j_URL:"
$ws.Range("A23").Value = "bolt://ncidb-q325-c.nci.nih.gov:7687"
$ws.Range("A24").Value = "User_name:"
$ws.Range("A25").Value = "neo4j"
$ws.Range("A26").Value = "PWD:"
$ws.Range("A27").Value = "icdcDBneo4j0"
$ws.Range("A28").Value = "Cypher:"

# A29 holds an explicit empty string (the blank cypher query that triggered
# the error above). A plain "" assignment clears the cell entirely, so force
# a text value via the leading-apostrophe trick and then strip the
# resulting quote-prefix formatting, leaving a true empty shared string.
$ws.Range("A29").Value = "'"
$ws.Range("A29").ClearFormats()

$ws.Range("A30").Value = "Output:"
$ws.Range("A31").Value = $outputPath
